# Insert a new weekly price record as row 759 in the "Papa" (potato)
# consolidated sheet, pushing every subsequent row down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 759:851 down to 760:852, carrying formatting (date style, etc.)
$ws.Rows(759).Insert()

# Populate the newly inserted row 759 with the new data point.
$ws.Range("A759").Value = 5
$ws.Range("B759").Value = "Macroferia Regional de Talca"
$ws.Range("C759").Value = "Maule"
$ws.Range("D759").Value = 45124
$ws.Range("E759").Value = 7
$ws.Range("F759").Value = 100114001
$ws.Range("G759").Value = "Papa"
$ws.Range("H759").Value = "Asterix"
$ws.Range("I759").Value = "1a (cosecha)"
$ws.Range("J759").Value = 1500
$ws.Range("K759").Value = 16000
$ws.Range("L759").Value = 16000
$ws.Range("M759").Value = 16000
$ws.Range("N759").Value = '$/saco 25 kilos'
$ws.Range("O759").Value = "Región del Maule"
$ws.Range("P759").Value = 640
$ws.Range("Q759").Value = 25
$ws.Range("R759").Value = "Hortaliza"
